# Fix: modified the title of excel data
# Renames the "Telefonnummer(+43...)" header (G1) to "Mobilnummer(+43...)"
# and moves the active selection to G4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header cell G1 from "Telefonnummer(+43…)" to "Mobilnummer(+43…)"
$ws.Range("G1").Value = "Mobilnummer(+43…)"

# Move the active selection to G4 (matches the diff's updated <selection> element)
$ws.Range("G4").Select()
